$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 10138
$ws.Range("I2").Value = 25045
$ws.Range("K2").Value = 25045
$ws.Range("M2").Value = -24932
$ws.Range("H18").Value = 1032.5
$ws.Range("I18").Value = 1032.5
$ws.Range("K18").Value = 1032.5
$ws.Range("M18").Value = -748.5
$ws.Range("H28").Value = 2789.4
$ws.Range("I28").Value = 1483.25
$ws.Range("J28").Value = 4282.143
$ws.Range("K28").Value = 1483.25
$ws.Range("L28").Value = 4282.143
$ws.Range("M28").Value = -998.25
$ws.Range("N28").Value = -5252.143
$ws.Range("H33").Value = 114.25
$ws.Range("I33").Value = 69
$ws.Range("K33").Value = 69
$ws.Range("M33").Value = 160
$ws.Range("H70").Value = 3923.1667
$ws.Range("I70").Value = 1752.5714
$ws.Range("K70").Value = 5257.7142
$ws.Range("M70").Value = -4987.7142
$ws.Range("H73").Value = 3923.1667
$ws.Range("I73").Value = 1752.5714
$ws.Range("K73").Value = 5257.7142
$ws.Range("M73").Value = -4321.7142
$ws.Range("H80").Value = 824.6667
$ws.Range("J80").Value = 912.25
$ws.Range("L80").Value = 2736.75
$ws.Range("N80").Value = -4732.75
$ws.Range("H83").Value = 824.6667
$ws.Range("J83").Value = 912.25
$ws.Range("L83").Value = 8210.25
$ws.Range("N83").Value = -18194.25
$ws.Range("H98").Value = 244.33333
$ws.Range("I98").Value = 244.33333
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 244.33333
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = $null
$ws.Range("N98").Value = 1253.66667
$ws.Range("H101").Value = 431.8
$ws.Range("J101").Value = 590
$ws.Range("L101").Value = 1770
$ws.Range("N101").Value = -5014
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = $null
$ws.Range("N105").Value = 0
$ws.Range("H122").Value = 244.33333
$ws.Range("I122").Value = 244.33333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 732.99999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = 1717.00001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 956.44446
$ws.Range("I2").Value = 964.25
$ws.Range("J2").Value = 894
$ws.Range("K2").Value = 964.25
$ws.Range("L2").Value = 894
$ws.Range("M2").Value = -851.25
$ws.Range("N2").Value = -1120
$ws.Range("H12").Value = 10000
$ws.Range("I12").Value = 10000
$ws.Range("K12").Value = 10000
$ws.Range("M12").Value = -9827
$ws.Range("H74").Value = 13100
$ws.Range("I74").Value = 18000
$ws.Range("J74").Value = 8200
$ws.Range("K74").Value = 18000
$ws.Range("L74").Value = 8200
$ws.Range("M74").Value = -17126
$ws.Range("N74").Value = -9948
$ws.Range("H77").Value = 13100
$ws.Range("I77").Value = 18000
$ws.Range("J77").Value = 8200
$ws.Range("K77").Value = 90000
$ws.Range("L77").Value = 41000
$ws.Range("M77").Value = -85632
$ws.Range("N77").Value = -49736
$ws.Range("H116").Value = 956.44446
$ws.Range("I116").Value = 964.25
$ws.Range("J116").Value = 894
$ws.Range("K116").Value = 964.25
$ws.Range("L116").Value = 894
$ws.Range("M116").Value = 1329.75
$ws.Range("N116").Value = -5482
$ws.Range("H132").Value = 2198.8
$ws.Range("I132").Value = 2298.5
$ws.Range("K132").Value = 6895.5
$ws.Range("M132").Value = -4365.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 956.44446
$ws.Range("I3").Value = 964.25
$ws.Range("J3").Value = 894
$ws.Range("K3").Value = 964.25
$ws.Range("L3").Value = 894
$ws.Range("M3").Value = -850.25
$ws.Range("N3").Value = -1122
$ws.Range("H86").Value = 3368.5715
$ws.Range("I86").Value = 2018.3334
$ws.Range("K86").Value = 2018.3334
$ws.Range("M86").Value = -895.3334
$ws.Range("H89").Value = 3368.5715
$ws.Range("I89").Value = 2018.3334
$ws.Range("K89").Value = 10091.667
$ws.Range("M89").Value = -4475.666999999999
$ws.Range("H94").Value = 584.8889
$ws.Range("I94").Value = 533.125
$ws.Range("K94").Value = 533.125
$ws.Range("M94").Value = -82.125

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 9418.362999999999
$ws.Range("J7").Value = 510
$ws.Range("L7").Value = 510
$ws.Range("N7").Value = -736
$ws.Range("H22").Value = 2887
$ws.Range("I22").Value = 1999
$ws.Range("J22").Value = 4219
$ws.Range("K22").Value = 1999
$ws.Range("L22").Value = 4219
$ws.Range("M22").Value = -1649
$ws.Range("N22").Value = -4919
$ws.Range("H107").Value = 756.46155
$ws.Range("I107").Value = 276.44446
$ws.Range("K107").Value = 276.44446
$ws.Range("M107").Value = 1643.55554

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 106.833336
$ws.Range("J2").Value = 72.40000000000001
$ws.Range("L2").Value = 434.4
$ws.Range("N2").Value = -660.4000000000001
$ws.Range("H4").Value = 166834160
$ws.Range("I4").Value = 200200600
$ws.Range("J4").Value = 2000
$ws.Range("K4").Value = 600601800
$ws.Range("L4").Value = 6000
$ws.Range("M4").Value = -600601688
$ws.Range("N4").Value = -6224
$ws.Range("H5").Value = 1125.8572
$ws.Range("I5").Value = 956.4
$ws.Range("J5").Value = 1549.5
$ws.Range("K5").Value = 2869.2
$ws.Range("L5").Value = 4648.5
$ws.Range("M5").Value = -2757.2
$ws.Range("N5").Value = -4872.5
$ws.Range("H12").Value = 175.5
$ws.Range("J12").Value = 191.2
$ws.Range("L12").Value = 573.5999999999999
$ws.Range("N12").Value = -919.5999999999999
$ws.Range("H38").Value = 399.6
$ws.Range("J38").Value = 335.6
$ws.Range("L38").Value = 1006.8
$ws.Range("N38").Value = -1700.8
$ws.Range("H86").Value = 1387.1428
$ws.Range("I86").Value = 1160
$ws.Range("K86").Value = 3480
$ws.Range("M86").Value = -2294
$ws.Range("H89").Value = 1387.1428
$ws.Range("I89").Value = 1160
$ws.Range("K89").Value = 10440
$ws.Range("M89").Value = -4512
$ws.Range("H112").Value = 3200
$ws.Range("J112").Value = 3200
$ws.Range("L112").Value = 9600
$ws.Range("N112").Value = -11816
$ws.Range("H114").Value = 265.16666
$ws.Range("I114").Value = 265.16666
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 795.4999799999999
$ws.Range("L114").Value = 0
$ws.Range("M114").Value = $null
$ws.Range("N114").Value = 2458.50002
$ws.Range("H135").Value = 1125.8572
$ws.Range("I135").Value = 956.4
$ws.Range("J135").Value = 1549.5
$ws.Range("K135").Value = 8607.6
$ws.Range("L135").Value = 13945.5
$ws.Range("M135").Value = -6072.6
$ws.Range("N135").Value = -19015.5
$ws.Range("H141").Value = 7000
$ws.Range("J141").Value = 12500
$ws.Range("L141").Value = 37500
$ws.Range("N141").Value = -47860

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 229.66667
$ws.Range("I2").Value = 170.7
$ws.Range("J2").Value = 303.375
$ws.Range("K2").Value = 170.7
$ws.Range("L2").Value = 303.375
$ws.Range("M2").Value = -57.69999999999999
$ws.Range("N2").Value = -529.375
$ws.Range("H80").Value = 3863.7144
$ws.Range("I80").Value = 3850
$ws.Range("K80").Value = 3850
$ws.Range("M80").Value = -2852
$ws.Range("H83").Value = 3863.7144
$ws.Range("I83").Value = 3850
$ws.Range("K83").Value = 19250
$ws.Range("M83").Value = -14258
$ws.Range("H107").Value = 3770.4285
$ws.Range("I107").Value = 3565.5
$ws.Range("J107").Value = 5000
$ws.Range("K107").Value = 3565.5
$ws.Range("L107").Value = 5000
$ws.Range("M107").Value = -1645.5
$ws.Range("N107").Value = -8840
$ws.Range("H122").Value = 5874.25
$ws.Range("I122").Value = 5398.8
$ws.Range("K122").Value = 16196.4
$ws.Range("M122").Value = -13746.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1162.5
$ws.Range("J22").Value = 1233.3334
$ws.Range("L22").Value = 1233.3334
$ws.Range("N22").Value = -1823.3334
$ws.Range("H27").Value = 1162.5
$ws.Range("J27").Value = 1233.3334
$ws.Range("L27").Value = 1233.3334
$ws.Range("N27").Value = -1447.3334
$ws.Range("H46").Value = 7091.6665
$ws.Range("I46").Value = 4183.3335
$ws.Range("J46").Value = 10000
$ws.Range("K46").Value = 4183.3335
$ws.Range("L46").Value = 10000
$ws.Range("M46").Value = -3995.3335
$ws.Range("N46").Value = -10376
$ws.Range("H55").Value = 917.6875
$ws.Range("I55").Value = 1065.6
$ws.Range("J55").Value = 671.1667
$ws.Range("K55").Value = 1065.6
$ws.Range("L55").Value = 671.1667
$ws.Range("M55").Value = -892.5999999999999
$ws.Range("N55").Value = -1017.1667
$ws.Range("H68").Value = 4931.25
$ws.Range("I68").Value = 1816.6666
$ws.Range("J68").Value = 6800
$ws.Range("K68").Value = 1816.6666
$ws.Range("L68").Value = 6800
$ws.Range("M68").Value = -1067.6666
$ws.Range("N68").Value = -8298
$ws.Range("H71").Value = 4931.25
$ws.Range("I71").Value = 1816.6666
$ws.Range("J71").Value = 6800
$ws.Range("K71").Value = 9083.333000000001
$ws.Range("L71").Value = 34000
$ws.Range("M71").Value = -5339.333000000001
$ws.Range("N71").Value = -41488

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1236.2632
$ws.Range("I132").Value = 1082.2667
$ws.Range("J132").Value = 1813.75
$ws.Range("K132").Value = 3246.800099999999
$ws.Range("L132").Value = 5441.25
$ws.Range("M132").Value = -716.8000999999995
$ws.Range("N132").Value = -10501.25
$ws.Range("H136").Value = 3585.8
$ws.Range("I136").Value = 2483.5
$ws.Range("J136").Value = 7995
$ws.Range("K136").Value = 7450.5
$ws.Range("L136").Value = 23985
$ws.Range("M136").Value = -4900.5
$ws.Range("N136").Value = -29085
